$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Guangming%Ye%NULL%2,                           Zhenyu%Pan%NULL%2,                           Yunbao%Pan%NULL%2,                           Qiaoling%Deng%NULL%2,                           Liangjun%Chen%NULL%2,                           Jin%Li%NULL%2,                           Yirong%Li%NULL%4,                           Xinghuan%Wang%NULL%5]"
$ws.Range("E3").Value = "[ Y.%Wang%null%1,                            W.% Jiang%null%1,                            Q.% He%null%1,                            C.% Wang%null%2,                            B.% Wang%null%2,                            P.% Zhou%null%1,                            N.% Dong%null%1,                            Q. % Tong%null%1,         Y.%Wang%null%1,         W.% Jiang%null%1,         Q.% He%null%1,         C.% Wang%null%2,         B.% Wang%null%2,         P.% Zhou%null%2,         N.% Dong%null%1,         Q. % Tong%null%1]"
$ws.Range("E4").Value = "[Yuanyuan%Xing%NULL%1,                           Pingzheng%Mo%NULL%2,                           Yu%Xiao%NULL%2,                           Oiu%Zhao%NULL%1,                           Yongxi%Zhang%NULL%4,                           Fan%Wang%NULL%2]"
$ws.Range("E5").Value = "[Kelvin Kai-Wang%To%NULL%2,                           Owen Tak-Yin%Tsang%NULL%1,                           Wai-Shing%Leung%NULL%1,                           Anthony Raymond%Tam%NULL%1,                           Tak-Chiu%Wu%NULL%1,                           David Christopher%Lung%NULL%1,                           Cyril Chik-Yan%Yip%NULL%1,                           Jian-Piao%Cai%NULL%2,                           Jacky Man-Chun%Chan%NULL%1,                           Thomas Shiu-Hong%Chik%NULL%1,                           Daphne Pui-Ling%Lau%NULL%1,                           Chris Yau-Chung%Choi%NULL%1,                           Lin-Lei%Chen%NULL%1,                           Wan-Mui%Chan%NULL%1,                           Kwok-Hung%Chan%NULL%2,                           Jonathan Daniel%Ip%NULL%1,                           Anthony Chin-Ki%Ng%NULL%1,                           Rosana Wing-Shan%Poon%NULL%1,                           Cui-Ting%Luo%NULL%1,                           Vincent Chi-Chung%Cheng%NULL%1,                           Jasper Fuk-Woo%Chan%NULL%3,                           Ivan Fan-Ngai%Hung%NULL%2,                           Zhiwei%Chen%NULL%1,                           Honglin%Chen%NULL%1,                           Kwok-Yung%Yuen%NULL%7]"
$ws.Range("E6").Value = "[ Bin%Ju%null%1,                  Qi%Zhang%null%1,                  Jiwan%Ge%null%1,                  Ruoke%Wang%null%1,                  Jing%Sun%null%1,                  Xiangyang%Ge%null%1,                  Jiazhen%Yu%null%1,                  Sisi%Shan%null%1,                  Bing%Zhou%null%1,                  Shuo%Song%null%1,                  Xian%Tang%null%1,                  Jinfang%Yu%null%1,                  Jun%Lan%null%1,                  Jing%Yuan%null%1,                  Haiyan%Wang%null%1,                  Juanjuan%Zhao%null%1,                  Shuye%Zhang%null%1,                  Youchun%Wang%null%1,                  Xuanling%Shi%null%1,                  Lei%Liu%null%1,                  Jincun%Zhao%null%1,                  Xinquan%Wang%null%1,                  Zheng%Zhang%null%1,                  Linqi%Zhang%null%1]"
$ws.Range("E7").Value = "[Chek Meng%Poh%NULL%1,                           Guillaume%Carissimo%NULL%1,                           Bei%Wang%NULL%2,                           Bei%Wang%NULL%0,                           Siti Naqiah%Amrun%NULL%1,                           Cheryl Yi-Pin%Lee%NULL%1,                           Rhonda Sin-Ling%Chee%NULL%1,                           Siew-Wai%Fong%NULL%1,                           Nicholas Kim-Wah%Yeo%NULL%1,                           Wen-Hsin%Lee%NULL%1,                           Anthony%Torres-Ruesta%NULL%1,                           Yee-Sin%Leo%NULL%2,                           Yee-Sin%Leo%NULL%0,                           Mark I-Cheng%Chen%NULL%1,                           Seow-Yen%Tan%NULL%1,                           Louis Yi Ann%Chai%NULL%1,                           Shirin%Kalimuddin%NULL%1,                           Shirley Seah Gek%Kheng%NULL%1,                           Siew-Yee%Thien%NULL%1,                           Barnaby Edward%Young%NULL%2,                           David C.%Lye%NULL%1,                           Brendon John%Hanson%NULL%1,                           Cheng-I%Wang%NULL%1,                           Laurent%Renia%NULL%1,                           Lisa F. P.%Ng%lisa_ng@immunol.a-star.edu.sg%2,                           Lisa F. P.%Ng%lisa_ng@immunol.a-star.edu.sg%0]"
$ws.Range("E8").Value = "[Wei%Zhang%NULL%1,                           Rong-Hui%Du%NULL%1,                           Bei%Li%NULL%1,                           Xiao-Shuang%Zheng%NULL%1,                           Xing-Lou%Yang%NULL%1,                           Ben%Hu%NULL%1,                           Yan-Yi%Wang%NULL%1,                           Geng-Fu%Xiao%NULL%1,                           Bing%Yan%NULL%1,                           Zheng-Li%Shi%NULL%1,                           Peng%Zhou%NULL%2,                           Peng%Zhou%NULL%0]"
$ws.Range("E9").Value = "[ Fatima%Amanat%null%1,                  Daniel%Stadlbauer%null%1,                  Shirin%Strohmeier%null%1,                  Thi H. O.%Nguyen%null%1,                  Veronika%Chromikova%null%1,                  Meagan%McMahon%null%1,                  Kaijun%Jiang%null%1,                  Guha Asthagiri%Arunkumar%null%1,                  Denise%Jurczyszak%null%1,                  Jose%Polanco%null%1,                  Maria%Bermudez-Gonzalez%null%1,                  Giulio%Kleiner%null%1,                  Teresa%Aydillo%null%1,                  Lisa%Miorin%null%1,                  Daniel S.%Fierer%null%1,                  Luz Amarilis%Lugo%null%1,                  Erna Milunka%Kojic%null%1,                  Jonathan%Stoever%null%1,                  Sean T. H.%Liu%null%1,                  Charlotte%Cunningham-Rundles%null%1,                  Philip L.%Felgner%null%1,                  Thomas%Moran%null%1,                  Adolfo%García-Sastre%null%1,                  Daniel%Caplivski%null%1,                  Allen C.%Cheng%null%1,                  Katherine%Kedzierska%null%1,                  Olli%Vapalahti%null%1,                  Jussi M.%Hepojoki%null%1,                  Viviana%Simon%null%1,                  Florian%Krammer%null%1]"
$ws.Range("E10").Value = "[Yicheng%Fang%NULL%3,                           Huangqi%Zhang%NULL%2,                          Huangqi%Zhang%NULL%0,                           Jicheng%Xie%NULL%2,                           Jicheng%Xie%NULL%0,                           Minjie%Lin%NULL%2,                           Minjie%Lin%NULL%0,                           Lingjun%Ying%NULL%2,                           Lingjun%Ying%NULL%0,                           Peipei%Pang%NULL%2,                           Peipei%Pang%NULL%0,                           Wenbin%Ji%1224190004@qq.com%2,                           Wenbin%Ji%1224190004@qq.com%0]"
